$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "MH0057:MH0059"
$ws.Range("E4").Value = "MH0054:MH0060"
$ws.Range("E5").Value = "MH0004:MH0005:MH0006:MH0007:MH0008:MH0012"
$ws.Range("E6").Value = "GH0012:GH0013:GH0009:GH0010"
$ws.Range("E7").Value = "GH0001:GH0002"
$ws.Range("F7").Value = "GH0003:GH0004:GH0005:GH0006:GH0007"
